# ADD results from server
# Update computed result values on sheets "2025", "2030", "2035"
$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1638.785132452183
$ws.Range("E2").Value = 293730.4047649619
$ws.Range("G2").Value = 80959.25712662081
$ws.Range("I2").Value = 142264.654638
$ws.Range("L2").Value = 525652.4050318201
$ws.Range("M2").Value = 111229.368349
$ws.Range("N2").Value = 70782.65578705262
$ws.Range("O2").Value = 67193.45619106332

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 34573.22838957706
$ws.Range("E2").Value = 169037.9833850653
$ws.Range("I2").Value = 133904.8587113994
$ws.Range("L2").Value = 152852.8481221266
$ws.Range("M2").Value = 60332.81188111824
$ws.Range("N2").Value = 21692.66563562064
$ws.Range("O2").Value = 11699.32491036274

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 19961.89159356778
$ws.Range("B2").Value = 18679.35888719985
$ws.Range("E2").Value = 121923.8325773797
$ws.Range("I2").Value = 170357.2063739901
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 56746.14559950977
$ws.Range("N2").Value = 44112.65282947898
$ws.Range("O2").Value = 52196.45092767161
